# feat：update water tower 1030
#
# Adds a new row (36) to the Tower table (Sheet1) describing water tower
# "1030" / 水龙男6, mirroring the layout/styles of the row above it (1029 /
# 水龙娘5), and moves the sheet's viewport/selection down to the newly
# added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the row immediately above the new data so the new row inherits
# the same per-column cell styles (s="1"/"9"/"15") as row 35.
$ws.Rows.Item(35).Copy()
$ws.Rows.Item(36).Insert()

# --- Row 36 values -------------------------------------------------------
$ws.Range("A36").Value = 1030
$ws.Range("B36").Value = "水龙男6"
$ws.Range("C36").Value = "水龙男6"
$ws.Range("D36").Value = 300
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = "1|2"
$ws.Range("G36").Value = 130042
$ws.Range("H36").Value = "Tower_attackTags_1|Tower_attackTags_2|Tower_attackTags_3|Tower_attackTags_4"
$ws.Range("I36").Value = "attackDamage|attackTime|attackCount|findRange"
$ws.Range("J36").Value = "300|400|520"
$ws.Range("K36").Value = "300|400|520"
$ws.Range("L36").Value = "A012FACE4430191FBD317598F7684D86|E718B09E4408CE5534779780E5365B64|E456238842ACC53D8C01EAABD11B256C"
$ws.Range("P36").Value = "F5DBBEBC4F82864C959DB8AFDD5CF39F"
$ws.Range("Q36").Value = 16
$ws.Range("T36").Value = 254276
$ws.Range("U36").Value = 1
$ws.Range("V36").Value = 0
$ws.Range("W36").Value = "3|3|3"
$ws.Range("X36").Value = "1|1|1"
$ws.Range("Y36").Value = "0|0|0"
$ws.Range("Z36").Value = "3|3|3"
$ws.Range("AA36").Value = "600|900|1300"
$ws.Range("AB36").Value = 20301
$ws.Range("AC36").Value = "FF9C6CFF|65AAFFFF|FF45FEFF"
$ws.Range("AD36").Value = 285705
$ws.Range("AE36").Value = 3
$ws.Range("AF36").Value = "1|1"

# --- Viewport: scroll/select down to the newly added row -----------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.DisplayHeadings = $true
$win.DisplayZeros = $true
$ws.Range("A22").Select()
$win.ScrollRow = 22
$win.ScrollColumn = 25
$ws.Range("AG44").Select()
